$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number:
# force Text format first, then restore the default "Normal" style
# after assigning the value so no visible style change remains.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.331.78"
$ws.Range("D3").Value = "2.036.62"
$ws.Range("E3").Value = "  +4.04%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "248.08"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("D7").Value = "60.67"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.398"
$ws.Range("E9").Value = "  +6.10%  "
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("E12").Value = "  +8.21%  "
$ws.Range("D13").Value = "0.864"
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").Value = "2.335.13"
$ws.Range("E14").Value = "  +3.86%  "
$ws.Range("D15").Value = "22.38"
$ws.Range("E15").Value = "  +3.39%  "
$ws.Range("D16").Value = "5.53"
$ws.Range("E16").Value = "  +5.56%  "
$ws.Range("D17").Value = "2.032.08"
$ws.Range("E17").Value = "  +3.75%  "
$ws.Range("D18").Value = "37.230.90"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").Value = "70.85"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").Value = "0.0₃0870"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").Value = "5.27"
$ws.Range("E21").Value = "  +3.86%  "
$ws.Range("D22").Value = "231.49"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D24").Value = "2.53"
$ws.Range("E24").Value = "  +4.29%  "
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").Value = "9.53"
$ws.Range("E26").Value = "  +3.83%  "
$ws.Range("D27").Value = "164.04"
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("D28").Value = "0.137"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("D29").Value = "19.91"
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("D30").Value = "1.39"
$ws.Range("E30").Value = "  +5.01%  "
$ws.Range("E31").Value = "  +2.82%  "
$ws.Range("D32").Value = "4.85"
$ws.Range("E32").Value = "  +2.14%  "
$ws.Range("D33").Value = "0.0672"
$ws.Range("E33").Value = "  +9.77%  "
$ws.Range("E34").Value = "  +11.61%  "
$ws.Range("D35").Value = "4.55"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("D36").Value = "3.59"
$ws.Range("E36").Value = "  +3.29%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "0.0984"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("D41").Value = "2.97"
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("D42").Value = "17.22"
$ws.Range("E42").Value = "  +8.94%  "
$ws.Range("E43").Value = "  +3.04%  "
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("D45").Value = "93.13"
$ws.Range("E45").Value = "  +4.96%  "
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("D47").Value = "1.390.99"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("E48").Value = "  +6.02%  "
$ws.Range("D49").Value = "2.17"
$ws.Range("E49").Value = "  +19.25%  "
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("D51").Value = "46.61"
$ws.Range("E51").Value = "  +1.95%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
